# Remove the seven rows that were dropped from the "Export" sheet.
# Deleting from the bottom up so earlier row numbers stay valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(27, 8, 7, 6, 5, 4, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).EntireRow.Delete()
}
